# v0.12: Add ice square background.
# Adds a new "ice" background-square tile to the "background" sheet (cell
# C2, next to the existing "bg" tile in B2) and makes "background" the
# active/selected sheet (it was "maze" before).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("background")

# New shared-string cell for the "ice" square, same look as the existing
# "bg" square in B2 (style carries over automatically since we just set
# the cell's value).
$ws.Range("C2").Value = "ice"

# Make "background" the active sheet/tab (was "maze"), and move its
# selection to C3.
$ws.Activate()
$ws.Range("C3").Select() | Out-Null
